{"js": "// The commit normalizes the \"Output report in <Folder>[...]: <file>\" lines\n// so every one of them reads \"Output report in <Folder>/analysis folder: <file>\",\n// and (as a side effect of touching the nearby \"See attached file[s] in\n// submittal.\" text) collapses a couple of previously-split runs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the four \"Output report in ...\" paragraphs by their exact text and\n// remember, for each, which folder/file pair it belongs to and which of the\n// two historical phrasings it uses:\n//   style \"colon\"  -> \"Output report in <Folder>: <file>\"\n//   style \"folder\" -> \"Output report in <Folder> folder: <file>\"\nconst targets = [];\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t === \"Output report in PyBank: financial_analysis.txt\") {\n    targets.push({ index: i, folder: \"PyBank\", file: \"financial_analysis.txt\", style: \"colon\" });\n  } else if (t === \"Output report in PyPoll: elections_results.txt\") {\n    targets.push({ index: i, folder: \"PyPoll\", file: \"elections_results.txt\", style: \"colon\" });\n  } else if (t === \"Output report in PyBank folder: financial_analysis.txt\") {\n    targets.push({ index: i, folder: \"PyBank\", file: \"financial_analysis.txt\", style: \"folder\" });\n  } else if (t === \"Output report in PyPoll folder: elections_results.txt\") {\n    targets.push({ index: i, folder: \"PyPoll\", file: \"elections_results.txt\", style: \"folder\" });\n  }\n}\n\nfor (const target of targets) {\n  const paragraph = items[target.index];\n\n  if (target.style === \"colon\") {\n    // \"Output report in <Folder>: <file>\" -> \"Output report in <Folder>/analysis folder: <file>\"\n    const colonRange = paragraph.search(\": \" + target.file, { matchCase: true });\n    colonRange.load(\"text\");\n    await context.sync();\n    colonRange.items[0].insertText(\"/analysis folder: \" + target.file, Word.InsertLocation.replace);\n    await context.sync();\n\n    // The same section also has a \"See attached file\" + \"s\" + \" in \" run\n    // split two paragraphs above the \"Output report\" line; re-set that text\n    // in place so it collapses into a single \"attached files in \" run.\n    const attachedParagraph = items[target.index - 2];\n    const attachedRange = attachedParagraph.search(\"attached files in \", { matchCase: true });\n    attachedRange.load(\"text\");\n    await context.sync();\n    attachedRange.items[0].insertText(\"attached files in \", Word.InsertLocation.replace);\n    await context.sync();\n  } else {\n    // \"Output report in <Folder> folder: <file>\" -> \"Output report in <Folder>/analysis folder: <file>\"\n    // Step 1: merge \"Output report\" + \" in \" into a single run.\n    const leadRange = paragraph.search(\"Output report in \", { matchCase: true });\n    leadRange.load(\"text\");\n    await context.sync();\n    leadRange.items[0].insertText(\"Output report in \", Word.InsertLocation.replace);\n    await context.sync();\n\n    // Step 2: insert \"/analysis\" right before \" folder: <file>\".\n    const folderRange = paragraph.search(\" folder: \" + target.file, { matchCase: true });\n    folderRange.load(\"text\");\n    await context.sync();\n    folderRange.items[0].insertText(\"/analysis folder: \" + target.file, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# The commit normalizes the \"Output report in <Folder>[...]: <file>\" lines\n# so every one of them reads \"Output report in <Folder>/analysis folder: <file>\".\n# Two of the four live in paragraphs that already read\n# \"Output report in <Folder>: <file>\" (just need \"/analysis folder\" spliced in\n# before the colon); the other two read\n# \"Output report in <Folder> folder: <file>\" (need \"Output report\"/\" in \"\n# joined back together and \"/analysis\" spliced in before \" folder\").\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Snapshot paragraph text (without the trailing paragraph mark) together with\n# each paragraph's index so we can re-fetch a fresh Range after earlier edits\n# shift character offsets further down the document.\n$paraInfo = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text.TrimEnd([char]13, [char]7)\n    $paraInfo += , @{ Index = $i; Text = $text }\n}\n\nfor ($i = 0; $i -lt $paraInfo.Length; $i++) {\n    $text = $paraInfo[$i].Text\n    $style = $null\n    $folder = $null\n    $file = $null\n\n    if ($text -eq \"Output report in PyBank: financial_analysis.txt\") {\n        $style = \"colon\"; $folder = \"PyBank\"; $file = \"financial_analysis.txt\"\n    } elseif ($text -eq \"Output report in PyPoll: elections_results.txt\") {\n        $style = \"colon\"; $folder = \"PyPoll\"; $file = \"elections_results.txt\"\n    } elseif ($text -eq \"Output report in PyBank folder: financial_analysis.txt\") {\n        $style = \"folder\"; $folder = \"PyBank\"; $file = \"financial_analysis.txt\"\n    } elseif ($text -eq \"Output report in PyPoll folder: elections_results.txt\") {\n        $style = \"folder\"; $folder = \"PyPoll\"; $file = \"elections_results.txt\"\n    }\n\n    if ($null -eq $style) { continue }\n\n    $paraIndex = $paraInfo[$i].Index\n    $outputPara = $d.Paragraphs.Item($paraIndex)\n\n    if ($style -eq \"colon\") {\n        # \"Output report in <Folder>: <file>\" -> \"Output report in <Folder>/analysis folder: <file>\"\n        $find = $outputPara.Range.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Execute(\": $file\", $false, $false, $false, $false, $false, $true, 1, $false, \"/analysis folder: $file\", 2) | Out-Null\n\n        # Two paragraphs above sits \"See attached file\" + \"s\" + \" in \" (three\n        # runs); re-set that same text so it collapses into \"attached files in \".\n        $attachedParaIndex = $paraIndex - 2\n        $attachedPara = $d.Paragraphs.Item($attachedParaIndex)\n        $find2 = $attachedPara.Range.Find\n        $find2.ClearFormatting()\n        $find2.Replacement.ClearFormatting()\n        $find2.Execute(\"attached files in \", $false, $false, $false, $false, $false, $true, 1, $false, \"attached files in \", 2) | Out-Null\n    } else {\n        # \"Output report in <Folder> folder: <file>\" -> \"Output report in <Folder>/analysis folder: <file>\"\n        # Step 1: merge \"Output report\" + \" in \" into one run.\n        $find = $outputPara.Range.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Execute(\"Output report in \", $false, $false, $false, $false, $false, $true, 1, $false, \"Output report in \", 2) | Out-Null\n\n        # Step 2: splice \"/analysis\" in right before \" folder: <file>\".\n        $find2 = $outputPara.Range.Find\n        $find2.ClearFormatting()\n        $find2.Replacement.ClearFormatting()\n        $find2.Execute(\" folder: $file\", $false, $false, $false, $false, $false, $true, 1, $false, \"/analysis folder: $file\", 2) | Out-Null\n    }\n}\n"}
